# Add a new column BB (the next forecast-vintage column) to the sheet,
# mirroring the layout of the existing BA column.
#
# - BB1 gets the new vintage date (and copies BA1's style/formatting).
# - BB3..BB18 duplicate the corresponding BA row's value (unchanged
#   historical YoY figures for those vintages).
# - BB19..BB21 get updated forecast values for the newest vintage.
# - Rows 2 and 22 have no data cell beyond column A, so BB is left empty
#   there, matching column BA.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy BA1's formatting (style index) onto BB1, then set its own date value.
$ws.Range("BA1").Copy($ws.Range("BB1"))
$ws.Range("BB1").Value = 45986

# Duplicate BA's YoY values into BB for the rows that remain unchanged.
# (Use Value2 for the read - Value's getter misbehaves in this host.)
$unchangedRows = 3..18
foreach ($r in $unchangedRows) {
    $baCell = $ws.Cells.Item($r, 53)   # column BA
    $bbCell = $ws.Cells.Item($r, 54)   # column BB
    $bbCell.Value = $baCell.Value2
}

# New forecast values for the most recent rows.
$ws.Range("BB19").Value = 2.043309689777173
$ws.Range("BB20").Value = 0.9040423720836799
$ws.Range("BB21").Value = 0.8372806497576768
